$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Rename sheets (uppercase + accent correction)
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Keep "PAINEIS DARQ" as the active/selected sheet (first tab)
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
